$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Locate the paragraph that currently ends in "...and finally" / " bring
#    him to the opposite bank..." (it is split into two runs around the
#    "_GoBack" bookmark). Remove the bookmark and merge the two runs into a
#    single run containing the full sentence.
# ---------------------------------------------------------------------------
$oldFullText = "There is only one choice the man could make on picking his first item, and that would be the parrot. By doing this, the parrot is safe from the cat and the bag of seed is safe from the parrot. The second item he transports could be either the cat or the bag of seed. Now that he has the parrot and the seed of cat on the opposite bank he faces his first problem. Either of these pairings comes with consequences. So in order to avoid leaving the bad pairing he would have to take the parrot back with him to the original bank. He then could bring either the bag of seed or cat, depending on what he chose as a second item over to the opposite bank. This would prevent any consequences. So that would leave the man to go retrieve the parrot and finally bring him to the opposite bank to finish his overall goal of having all 3 items on the opposite bank."

# The "_GoBack" bookmark sits between the two runs; remove it first so that
# the paragraph text below is contiguous.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Re-typing the paragraph's full (now bookmark-free) text via Find/Replace
# collapses the two adjacent, identically formatted runs into one clean run.
$d.Content.Find.Execute($oldFullText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $oldFullText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Find the paragraph we just normalized, and insert a new paragraph right
#    after it (Word copies paragraph/run formatting automatically).
# ---------------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq $oldFullText) {
        $target = $p
    }
}

$afterRange = $target.Range
$afterRange.Collapse(0)
$afterRange.InsertParagraphAfter()

# The newly inserted paragraph is the one immediately following $target.
$newPara = $target.Next()

# ---------------------------------------------------------------------------
# 3. Fill the new paragraph with the test-theory text, using a temporary
#    trailing placeholder character so that the "_GoBack" bookmark can later
#    be anchored as a true zero-length range right before the paragraph
#    mark (inserting a bookmark directly at paragraph-end can otherwise land
#    incorrectly in this environment).
# ---------------------------------------------------------------------------
$newText = "In order to test this theory, I had to draw out the situation at hand. This would include the faulty pairings and the fact he could only bring one item with him at a time. When not bringing an item back with him the second trip, no matter what two items he brought over there would be a consequence leaving them alone. The solution only works if the parrot is brought over first, if not the parrot would be left with the bag of seed or the cat."

$placeholder = "\u0001"
$newParaRange = $newPara.Range
$newParaRange.MoveEnd(1, -1)
$newParaRange.Text = $newText + $placeholder

# ---------------------------------------------------------------------------
# 4. Insert the "_GoBack" bookmark immediately before the placeholder
#    character (a genuine, non-degenerate position), then delete the
#    placeholder so the bookmark ends up as a zero-length anchor right at
#    the end of the paragraph's text, matching the original document.
# ---------------------------------------------------------------------------
$newPara = $target.Next()
$bmPos = $newPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$bmRange.Bookmarks.Add("_GoBack")

$newPara = $target.Next()
$placeholderRange = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$placeholderRange.Delete()
